# Update "Forecast Comparison" sheet with corrected forecast output:
#  - insert a new "Week_Start_Date" column after "Week" (new column B)
#  - shorten week labels from "W01".."W16" to "W1".."W16"
#  - turn "is_holiday_week" into a real boolean column

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# 1. Insert a new column B for the week start date, shifting ASIN.. etc right.
$ws.Range("B1").EntireColumn.Insert()

# 2. Header for the new column.
$ws.Range("B1").Value = "Week_Start_Date"

# Force the new column's data cells to be plain text so ISO-looking dates
# ("2025-01-05") are kept as literal strings instead of being parsed into
# date serials.
$ws.Range("B2:B17").NumberFormat = "@"

$weekStarts = @(
    "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
    "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
    "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
    "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)

$weekLabels = @(
    "W1", "W2", "W3", "W4", "W5", "W6", "W7", "W8",
    "W9", "W10", "W11", "W12", "W13", "W14", "W15", "W16"
)

for ($i = 0; $i -lt $weekStarts.Length; $i++) {
    $row = $i + 2
    # 3. Shorten the week label in column A ("W01" -> "W1", etc.)
    $ws.Cells.Item($row, 1).Value = $weekLabels[$i]
    # 4. Populate the new Week_Start_Date column.
    $ws.Cells.Item($row, 2).Value = $weekStarts[$i]
    # 5. Re-type is_holiday_week (now column J) as a boolean, same values.
    $ws.Cells.Item($row, 10).Value = $false
}

"Done"
